$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with two new columns (P, Q), copying the header style (bold/border/center)
# from the existing last header cell (O1) then setting the new values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Fill the new data columns P and Q (rows 2-25) with value 2.
$ws.Range("P2:Q25").Value = 2

# Swap the values in columns I, K, M, O for data rows 2-25:
#   I: 1 -> 2
#   K: 2 -> 1
#   M: 1 -> 2
#   O: 2 -> 1
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1
